# Update "new Madigan bike hours" - Riders (C) and Average (D) columns
# for rows 2-28 on the "Ridership" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ridership")

$newValues = @{
    2  = @(200, 195.29)
    3  = @(196, 222.46)
    4  = @(257, 238.64)
    5  = @(107, 120.2)
    6  = @(43, 101.08)
    7  = @(175, 230.71)
    8  = @(131, 206.27)
    9  = @(260, 199.6)
    10 = @(271, 225.93)
    11 = @(232, 238.2)
    12 = @(95, 118.62)
    13 = @(127, 102.93)
    14 = @(208, 229.2)
    15 = @(242, 208.5)
    16 = @(264, 203.62)
    17 = @(282, 229.67)
    18 = @(246, 238.69)
    19 = @(143, 120.06)
    20 = @(96, 102.47)
    21 = @(219, 209.12)
    22 = @(288, 208.59)
    23 = @(241, 230.38)
    24 = @(258, 239.82)
    25 = @(98, 118.83)
    26 = @(93, 101.88)
    27 = @(242, 230)
    28 = @(238, 210.72)
}

foreach ($row in $newValues.Keys) {
    $vals = $newValues[$row]
    $ws.Cells.Item($row, 3).Value = $vals[0]
    $ws.Cells.Item($row, 4).Value = $vals[1]
}

$wb.Save()
